# Remove column M from the alcohol measurement data sheet.
# This deletes the entire column M (shifting column N, the last
# populated column, left into its place) and updates the active
# selection to reflect the new rightmost data column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("M:M").Delete()

$ws.Range("M1").Select() | Out-Null
